$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the easting/northing coordinates to whole numbers
$ws.Range("Q2").Value = 500898
$ws.Range("R2").Value = 6544336

# Clear the start/end time cells (Z2, AB2) entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
